$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 2271.5
$ws.Range("I34").Value = 2271.5
$ws.Range("K34").Value = 2271.5
$ws.Range("M34").Value = -2068.5
$ws.Range("H36").Value = 2271.5
$ws.Range("I36").Value = 2271.5
$ws.Range("K36").Value = 2271.5
$ws.Range("M36").Value = -1556.5
$ws.Range("H88").Value = 1472.0769
$ws.Range("I88").Value = 710.6
$ws.Range("J88").Value = 1948
$ws.Range("K88").Value = 710.6
$ws.Range("L88").Value = 1948
$ws.Range("M88").Value = -304.6
$ws.Range("N88").Value = -2760
$ws.Range("H91").Value = 1472.0769
$ws.Range("I91").Value = 710.6
$ws.Range("J91").Value = 1948
$ws.Range("K91").Value = 710.6
$ws.Range("L91").Value = 1948
$ws.Range("M91").Value = 693.4
$ws.Range("N91").Value = -4756
$ws.Range("H116").Value = 3408.5833
$ws.Range("I116").Value = 3674.625
$ws.Range("J116").Value = 2876.5
$ws.Range("K116").Value = 3674.625
$ws.Range("L116").Value = 2876.5
$ws.Range("M116").Value = -232.625
$ws.Range("N116").Value = -9760.5
$ws.Range("H125").Value = 2354.3635
$ws.Range("I125").Value = 1706.3077
$ws.Range("J125").Value = 3290.4443
$ws.Range("K125").Value = 15356.7693
$ws.Range("L125").Value = 29613.9987
$ws.Range("M125").Value = -12896.7693
$ws.Range("N125").Value = -34533.9987
$ws.Range("H132").Value = 336425.97
$ws.Range("I132").Value = 360385.1
$ws.Range("J132").Value = 998
$ws.Range("K132").Value = 1081155.3
$ws.Range("L132").Value = 2994
$ws.Range("M132").Value = -1078625.3
$ws.Range("N132").Value = -8054
$ws.Range("H141").Value = 1455.9584
$ws.Range("I141").Value = 1159.1428
$ws.Range("J141").Value = 1871.5
$ws.Range("K141").Value = 3477.4284
$ws.Range("L141").Value = 5614.5
$ws.Range("M141").Value = 1702.5716
$ws.Range("N141").Value = -15974.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 20666.334
$ws.Range("I21").Value = 15999
$ws.Range("J21").Value = 23000
$ws.Range("K21").Value = 15999
$ws.Range("L21").Value = 23000
$ws.Range("M21").Value = -15625
$ws.Range("N21").Value = -23748
$ws.Range("H29").Value = 2313.3333
$ws.Range("J29").Value = 2313.3333
$ws.Range("L29").Value = 2313.3333
$ws.Range("N29").Value = -2929.3333
$ws.Range("I30").Value = 5000
$ws.Range("J30").Value = 2000
$ws.Range("K30").Value = 5000
$ws.Range("L30").Value = 2000
$ws.Range("M30").Value = -4850
$ws.Range("N30").Value = -2300
$ws.Range("H61").Value = 13335191
$ws.Range("I61").Value = 16668429
$ws.Range("J61").Value = 2239.6
$ws.Range("K61").Value = 16668429
$ws.Range("L61").Value = 2239.6
$ws.Range("M61").Value = -16668217
$ws.Range("N61").Value = -2663.6
$ws.Range("H136").Value = 13335191
$ws.Range("I136").Value = 16668429
$ws.Range("J136").Value = 2239.6
$ws.Range("K136").Value = 50005287
$ws.Range("L136").Value = 6718.799999999999
$ws.Range("M136").Value = -50002737
$ws.Range("N136").Value = -11818.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 32542.4
$ws.Range("I134").Value = 45076
$ws.Range("K134").Value = 135228
$ws.Range("M134").Value = -132693

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1666.6072
$ws.Range("I132").Value = 1504.9788
$ws.Range("J132").Value = 2510.6667
$ws.Range("K132").Value = 4514.936400000001
$ws.Range("L132").Value = 7532.000100000001
$ws.Range("M132").Value = -1984.936400000001
$ws.Range("N132").Value = -12592.0001
$ws.Range("H134").Value = 1714.7667
$ws.Range("I134").Value = 1892.9131
$ws.Range("J134").Value = 1129.4286
$ws.Range("K134").Value = 5678.7393
$ws.Range("L134").Value = 3388.2858
$ws.Range("M134").Value = -3143.7393
$ws.Range("N134").Value = -8458.2858

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 1118.1666
$ws.Range("I26").Value = 77.25
$ws.Range("J26").Value = 3200
$ws.Range("K26").Value = 231.75
$ws.Range("L26").Value = 9600
$ws.Range("M26").Value = 56.25
$ws.Range("N26").Value = -10176
$ws.Range("H33").Value = 27778000
$ws.Range("I33").Value = 41666900
$ws.Range("K33").Value = 250001400
$ws.Range("M33").Value = -250001117
$ws.Range("H86").Value = 1600
$ws.Range("I86").Value = 300
$ws.Range("J86").Value = 2250
$ws.Range("K86").Value = 900
$ws.Range("L86").Value = 6750
$ws.Range("M86").Value = 286
$ws.Range("N86").Value = -9122
$ws.Range("H89").Value = 1600
$ws.Range("I89").Value = 300
$ws.Range("J89").Value = 2250
$ws.Range("K89").Value = 2700
$ws.Range("L89").Value = 20250
$ws.Range("M89").Value = 3228
$ws.Range("N89").Value = -32106
$ws.Range("H113").Value = 549.4091
$ws.Range("I113").Value = 598.75
$ws.Range("J113").Value = 521.2143
$ws.Range("K113").Value = 1796.25
$ws.Range("L113").Value = 1563.6429
$ws.Range("M113").Value = 373.75
$ws.Range("N113").Value = -5903.6429

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1612.2034
$ws.Range("I132").Value = 1661.5135
$ws.Range("J132").Value = 1529.2727
$ws.Range("K132").Value = 4984.5405
$ws.Range("L132").Value = 4587.8181
$ws.Range("M132").Value = -2454.5405
$ws.Range("N132").Value = -9647.8181

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 3200
$ws.Range("I4").Value = 3200
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 3200
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -3087
$ws.Range("N4").ClearContents()
$ws.Range("H5").Value = 3000
$ws.Range("J5").Value = 3000
$ws.Range("L5").Value = 3000
$ws.Range("N5").Value = -3226
$ws.Range("H7").Value = 1746.2
$ws.Range("I7").Value = 1811.3572
$ws.Range("J7").Value = 1594.1666
$ws.Range("K7").Value = 1811.3572
$ws.Range("L7").Value = 1594.1666
$ws.Range("M7").Value = -1699.3572
$ws.Range("N7").Value = -1818.1666
$ws.Range("H16").Value = 1906.9615
$ws.Range("I16").Value = 2013.375
$ws.Range("K16").Value = 2013.375
$ws.Range("M16").Value = -1843.375
$ws.Range("H26").Value = 14998.667
$ws.Range("J26").Value = 14998.667
$ws.Range("L26").Value = 14998.667
$ws.Range("N26").Value = -15588.667
$ws.Range("H28").Value = 3200
$ws.Range("I28").Value = 3200
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 3200
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -2968
$ws.Range("N28").ClearContents()
$ws.Range("H31").Value = 3207.5
$ws.Range("I31").Value = 765
$ws.Range("J31").Value = 5650
$ws.Range("K31").Value = 765
$ws.Range("L31").Value = 5650
$ws.Range("M31").Value = -517
$ws.Range("N31").Value = -6146
$ws.Range("H37").Value = 3200
$ws.Range("I37").Value = 3200
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 3200
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -3093
$ws.Range("N37").ClearContents()
$ws.Range("H126").Value = 1746.2
$ws.Range("I126").Value = 1811.3572
$ws.Range("J126").Value = 1594.1666
$ws.Range("K126").Value = 5434.071599999999
$ws.Range("L126").Value = 4782.4998
$ws.Range("M126").Value = -2964.071599999999
$ws.Range("N126").Value = -9722.4998
$ws.Range("H136").Value = 3260.262
$ws.Range("J136").Value = 3099.6
$ws.Range("L136").Value = 9298.799999999999
$ws.Range("N136").Value = -14398.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 70017
$ws.Range("J21").Value = 70017
$ws.Range("L21").Value = 70017
$ws.Range("N21").Value = -70487
$ws.Range("H24").Value = 70010
$ws.Range("J24").Value = 70010
$ws.Range("L24").Value = 70010
$ws.Range("N24").Value = -70470
$ws.Range("H29").Value = 63007.332
$ws.Range("J29").Value = 63007.332
$ws.Range("L29").Value = 63007.332
$ws.Range("N29").Value = -63587.332
$ws.Range("H35").Value = 70017
$ws.Range("J35").Value = 70017
$ws.Range("L35").Value = 70017
$ws.Range("N35").Value = -70597
$ws.Range("H132").Value = 7291.8
$ws.Range("I132").Value = 8090.381
$ws.Range("K132").Value = 24271.143
$ws.Range("M132").Value = -21741.143

Write-Host "Applied all changes"